# Rename the "AddressBook" class box to "Concierge" in the Logic Component
# Class Diagram (shape "Rectangle 62", id=16, on slide 1). The box has two
# paragraphs - "AddressBook" and "Parser" - and only the first paragraph's
# text changes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Find the shape whose first paragraph currently reads "AddressBook"
# (Paragraph.Text includes the trailing paragraph-mark character, so trim it
# before comparing).
$target = $null
foreach ($sh in $s.Shapes) {
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $firstPara = $sh.TextFrame.TextRange.Paragraphs(1).Text
        if ($firstPara.TrimEnd() -eq "AddressBook") {
            $target = $sh
        }
    }
}

$tr = $target.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)

# Delete the first paragraph ("AddressBook") outright and retype the new
# name as a brand-new run, rather than just reassigning .Text - this drops
# the stale spell-check flag (err="1") that PowerPoint had cached for the
# old word, matching how a fresh, correctly-spelled run is produced, and
# also drops the separate endParaRPr that isn't needed once the run itself
# carries through to the paragraph mark.
[void]$para1.Delete()
[void]$tr.InsertBefore("Concierge" + [char]13)

Write-Output $tr.Text
